# Add the new "Player Info" sheet before the existing "ODI Batting" sheet.
$wb = $excel.ActiveWorkbook
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "Player Info"

# Re-fetch both sheets by name: the COM handles captured before `Add()` track
# worksheet *position*, not identity, and `Add()` shifts everything right by one.
$playerInfo = $wb.Worksheets.Item("Player Info")
$odiBatting = $wb.Worksheets.Item("ODI Batting")

# Header row.
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold/centered/bordered header style already used on "ODI Batting" row 1.
$odiBatting.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122) # xlPasteFormats

# Data row.
$playerInfo.Range("A2").Formula = "'5845"
$playerInfo.Range("B2").Value = "Haider Ali"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"

# Strip the quote-prefix formatting iron_native applied for the numeric-looking
# literal in A2 so it falls back to the sheet's plain/default style, matching a
# freshly authored inline string cell.
$playerInfo.Range("B2").Copy()
$playerInfo.Range("A2").PasteSpecial(-4122) # xlPasteFormats

# Rework the "ODI Batting" sheet's MATCH_CARD_LINK column into a bare MATCH_CODE.
$odiBatting.Range("D1").Value = "MATCH_CODE"
$odiBatting.Range("D2").Formula = "'4433"
$odiBatting.Range("D3").Formula = "'4434"

$odiBatting.Range("E2").Copy()
$odiBatting.Range("D2").PasteSpecial(-4122) # xlPasteFormats
$odiBatting.Range("E3").Copy()
$odiBatting.Range("D3").PasteSpecial(-4122) # xlPasteFormats
